# Actualización SmartScore desde Streamlit (remas ali almadani)
# Adds a new response row (row 27) to the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

# --- Plain text / identity columns -----------------------------------
$ws.Range("A$row").Value = "remas ali almadani_20251202_134135"
$ws.Range("C$row").Value = "remas ali almadani"
$ws.Range("D$row").Value = 19
$ws.Range("E$row").Value = "Female"
$ws.Range("F$row").Value = "2025-12-02 13:41:35"

# --- B27: present but blank (no "Sin SmartScore" marker this time) ---
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("B$row").Value = ""
$ws.Range("B$row").Style = "Normal"

# --- G27: raw JSON weights blob used by the Streamlit app ------------
$weights = @"
{
  "portion": 0.2,
  "diet": 0.2857142857142857,
  "salt": 0.6,
  "fat": 0.6,
  "natural": 0.8,
  "convenience": 0.4,
  "price": 1.0
}
"@
$ws.Range("G$row").Value = $weights

# --- Instant Noodles -----------------------------------------------------
$ws.Range("H$row").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("J$row").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K$row").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("M$row").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("N$row").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("P$row").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

# --- Mac & Cheese ----------------------------------------------------
$ws.Range("Q$row").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("S$row").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("T$row").Value = "Annie’s Shells & White Cheddar"
$ws.Range("V$row").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("W$row").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("Y$row").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

# --- Ready to Eat ------------------------------------------------------
$ws.Range("Z$row").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AB$row").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC$row").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AE$row").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"
$ws.Range("AF$row").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AH$row").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

# --- SmartScore numeric-looking columns, stored as TEXT (matches the
#     export's inline-string cells, not real numbers) ------------------
$scoreCols = "I", "L", "O", "R", "U", "X", "AA", "AD", "AG"
$scoreVals = "0.578", "0.566", "0.455", "0.712", "0.625", "0.567", "0.657", "0.656", "0.644"

for ($i = 0; $i -lt $scoreCols.Length; $i++) {
    $addr = "$($scoreCols[$i])$row"
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $scoreVals[$i]
    $ws.Range($addr).Style = "Normal"
}

# The multi-line JSON in G27 makes the engine auto-grow the row height;
# AutoFit it back down so the row stays at the sheet's default height
# (matching every other data row, none of which carry an explicit height).
$ws.Rows.Item($row).AutoFit()

